# edit.ps1 - applies the "LnCo as on Dec 27" diff to Config.xlsx
$wb = $excel.ActiveWorkbook

# -------------------------------------------------------------------
# 1) "Main" sheet: populate new keys MB51_Movement_types_list (row 59)
#    and Inventory_mapping_exception_percentage (row 60)
# -------------------------------------------------------------------
$main = $wb.Worksheets.Item("Main")

$main.Range("A59").Value = "MB51_Movement_types_list"
$main.Range("B59").Value = "[101, 102, 122, 123]"

$main.Range("A60").Value = "Inventory_mapping_exception_percentage"
$main.Range("B60").Value = 10

# Update the sheet view's selection / scroll position
$main.Activate()
try {
    $excel.ActiveWindow.TopLeftCell = $main.Range("A50")
} catch {
}
$main.Range("B60").Select()

# -------------------------------------------------------------------
# 2) "Inventory Mapping" sheet: insert two new rows (15 & 16) holding
#    the "Movement Type Subject" / "Movement Type Body" mail texts,
#    pushing everything from the old row 16 onward down by two rows.
# -------------------------------------------------------------------
$inv = $wb.Worksheets.Item("Inventory Mapping")

$inv.Rows("16:17").Insert()

$inv.Range("A15").Value = "Movement Type Subject"
$inv.Range("B15").Value = "Inventory Mapping input file column data is empty"

$inv.Range("A16").Value = "Movement Type Body"
$inv.Range("B16").Value = "Hello,`nPurchase Register Movement Type Column data is empty. `nThanks & Regards,`nL & Co   `n                                                                                                                                                                                                                                                                                                                                                                                                                                                                                                                                                                                                                                                                                                                                                                                                           "

# Match formatting used by the equivalent Subject/Body pairs already
# present on the sheet (e.g. B9/B10, the "GR Document_Number" pair)
$inv.Range("B9").Copy()
$inv.Range("B15").PasteSpecial(-4122)
$inv.Range("B10").Copy()
$inv.Range("B16").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$inv.Rows(15).RowHeight = 100.8
$inv.Rows(16).RowHeight = 100.8

# Update the sheet view's selection
$inv.Activate()
$inv.Range("B8").Select()

$main.Activate()
